$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.380719000000001
$ws.Range("H2").Value = 28.142157
$ws.Range("I2").Value = 0.03679977590837273
$ws.Range("J2").Value = 0.03679977590837273
$ws.Range("M2").Value = 0.08962966666666666
$ws.Range("N2").Value = 0.268889
$ws.Range("O2").Value = 0.4339761198462219
$ws.Range("P2").Value = 0.4339761198462219
$ws.Range("Q2").Value = 0.8407907170636667
$ws.Range("R2").Value = 7.567116453573
$ws.Range("S2").Value = 0.01597022395992607
$ws.Range("T2").Value = 0.01597022395992607

$ws.Range("G3").Value = 9.380719000000001
$ws.Range("H3").Value = 28.142157
$ws.Range("I3").Value = 0.03679977590837273
$ws.Range("J3").Value = 0.03679977590837273
$ws.Range("M3").Value = 0.1169016666666667
$ws.Range("N3").Value = 0.350705
$ws.Range("O3").Value = 0.5660238801537781
$ws.Range("P3").Value = 0.5660238801537781
$ws.Range("Q3").Value = 1.096621685631667
$ws.Range("R3").Value = 9.869595170685
$ws.Range("S3").Value = 0.02082955194844666
$ws.Range("T3").Value = 0.02082955194844666

$ws.Range("I4").Value = 0.3547860986448385
$ws.Range("J4").Value = 0.3547860986448385
$ws.Range("M4").Value = 0.08962966666666666
$ws.Range("N4").Value = 0.268889
$ws.Range("O4").Value = 0.4339761198462219
$ws.Range("P4").Value = 0.4339761198462219
$ws.Range("Q4").Value = 8.106050945162
$ws.Range("R4").Value = 72.954458506458
$ws.Range("S4").Value = 0.1539686944652659
$ws.Range("T4").Value = 0.1539686944652659

$ws.Range("I5").Value = 0.3547860986448385
$ws.Range("J5").Value = 0.3547860986448385
$ws.Range("M5").Value = 0.1169016666666667
$ws.Range("N5").Value = 0.350705
$ws.Range("O5").Value = 0.5660238801537781
$ws.Range("P5").Value = 0.5660238801537781
$ws.Range("Q5").Value = 10.57251355289
$ws.Range("R5").Value = 95.15262197601
$ws.Range("S5").Value = 0.2008174041795726
$ws.Range("T5").Value = 0.2008174041795726

$ws.Range("G6").Value = 100.179423
$ws.Range("H6").Value = 300.538269
$ws.Range("I6").Value = 0.3929954960840508
$ws.Range("J6").Value = 0.3929954960840508
$ws.Range("M6").Value = 0.08962966666666666
$ws.Range("N6").Value = 0.268889
$ws.Range("O6").Value = 0.4339761198462219
$ws.Range("P6").Value = 0.4339761198462219
$ws.Range("Q6").Value = 8.979048290348999
$ws.Range("R6").Value = 80.811434613141
$ws.Range("S6").Value = 0.1705506605075974
$ws.Range("T6").Value = 0.1705506605075974

$ws.Range("G7").Value = 100.179423
$ws.Range("H7").Value = 300.538269
$ws.Range("I7").Value = 0.3929954960840508
$ws.Range("J7").Value = 0.3929954960840508
$ws.Range("M7").Value = 0.1169016666666667
$ws.Range("N7").Value = 0.350705
$ws.Range("O7").Value = 0.5660238801537781
$ws.Range("P7").Value = 0.5660238801537781
$ws.Range("Q7").Value = 11.711141514405
$ws.Range("R7").Value = 105.400273629645
$ws.Range("S7").Value = 0.2224448355764533
$ws.Range("T7").Value = 0.2224448355764533

$ws.Range("G8").Value = 1.427630666666667
$ws.Range("H8").Value = 4.282892
$ws.Range("I8").Value = 0.005600475679236752
$ws.Range("J8").Value = 0.005600475679236752
$ws.Range("M8").Value = 0.08962966666666666
$ws.Range("N8").Value = 0.268889
$ws.Range("O8").Value = 0.4339761198462219
$ws.Range("P8").Value = 0.4339761198462219
$ws.Range("Q8").Value = 0.1279580607764444
$ws.Range("R8").Value = 1.151622546988
$ws.Range("S8").Value = 0.002430472704568299
$ws.Range("T8").Value = 0.002430472704568299

$ws.Range("G9").Value = 1.427630666666667
$ws.Range("H9").Value = 4.282892
$ws.Range("I9").Value = 0.005600475679236752
$ws.Range("J9").Value = 0.005600475679236752
$ws.Range("M9").Value = 0.1169016666666667
$ws.Range("N9").Value = 0.350705
$ws.Range("O9").Value = 0.5660238801537781
$ws.Range("P9").Value = 0.5660238801537781
$ws.Range("Q9").Value = 0.1668924043177778
$ws.Range("R9").Value = 1.50203163886
$ws.Range("S9").Value = 0.003170002974668452
$ws.Range("T9").Value = 0.003170002974668452

$ws.Range("G10").Value = 53.48524799999999
$ws.Range("H10").Value = 160.455744
$ws.Range("I10").Value = 0.2098181536835013
$ws.Range("J10").Value = 0.2098181536835013
$ws.Range("M10").Value = 0.08962966666666666
$ws.Range("N10").Value = 0.268889
$ws.Range("O10").Value = 0.4339761198462219
$ws.Range("P10").Value = 0.4339761198462219
$ws.Range("Q10").Value = 4.793864949823999
$ws.Range("R10").Value = 43.14478454841599
$ws.Range("S10").Value = 0.09105606820886414
$ws.Range("T10").Value = 0.09105606820886414

$ws.Range("G11").Value = 53.48524799999999
$ws.Range("H11").Value = 160.455744
$ws.Range("I11").Value = 0.2098181536835013
$ws.Range("J11").Value = 0.2098181536835013
$ws.Range("M11").Value = 0.1169016666666667
$ws.Range("N11").Value = 0.350705
$ws.Range("O11").Value = 0.5660238801537781
$ws.Range("P11").Value = 0.5660238801537781
$ws.Range("Q11").Value = 6.252514633279999
$ws.Range("R11").Value = 56.27263169951999
$ws.Range("S11").Value = 0.1187620854746371
$ws.Range("T11").Value = 0.1187620854746371
